# Generate Report for Handback
# Updates the "eee1ae94-e31a-4bfe-92c1-8485cb9fbde4.md" row (row 5) across the
# Overview / zh-cn / de-de sheets to reflect a completed handback:
#   - Status columns move from "Ready for handoff" to
#     "Handed back: in sync with en-US"
#   - Latest Handback DateTime is refreshed
#   - Error Detail (stale-handback warning) is cleared

$wb = $excel.ActiveWorkbook

$status = "Handed back: in sync with en-US"

# --- Overview sheet ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E5").Value = $status
$overview.Range("F5").Value = $status

# --- zh-cn sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C5").Value = $status
$zhcn.Range("K5").Value = "2016-09-07 09:13:26"
$zhcn.Range("P5").Value = ""

# --- de-de sheet ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C5").Value = $status
$dede.Range("K5").Value = "2016-09-07 09:13:45"
$dede.Range("P5").Value = ""
